# The sheet's columns D,E,F,G ("codeforiati:category-name",
# "codeforiati:category-code", "codeforiati:group-name",
# "codeforiati:group-code" plus their data) are right-rotated by one
# column for every row (including the header): the old G value becomes
# the new D value, and D,E,F each shift one column to the right
# (D->E, E->F, F->G).
#
# We use Range.Copy(Destination) (rather than .Value/.Value2 assignment)
# so that numeric-looking text such as "110" or "111" stays a text
# (shared-string) cell instead of being auto-coerced into a number, and
# so no cell style/number-format gets introduced along the way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

$srcRangeG = $ws.Range("G1:G$lastRow")
$tempRangeI = $ws.Range("I1:I$lastRow")
$srcRangeDF = $ws.Range("D1:F$lastRow")
$dstRangeEG = $ws.Range("E1:G$lastRow")
$dstRangeD = $ws.Range("D1:D$lastRow")

# 1) Stash the original column G (group-code) in scratch column I.
$srcRangeG.Copy($tempRangeI)

# 2) Shift D:F right into E:G (category-name/category-code/group-name
#    each move one column over).
$srcRangeDF.Copy($dstRangeEG)

# 3) Drop the stashed original G values into D (group-code becomes the
#    new first of the four columns).
$tempRangeI.Copy($dstRangeD)

# 4) Clean up the scratch column.
$tempRangeI.Clear()
